$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Activate()
$ws2.Range("A4").Value = 5
$ws2.Range("B4").Formula = "=A4*7"
$ws2.Range("A5").Select()

$ws1.Activate()
$ws1.Range("A1").Value = "Hello"
$ws1.Range("A2").Value = "World"
$ws1.Range("A3").Value = "!"
$ws1.Range("A4").Select()
